$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1430.2941
$ws.Range("J112").Value = 1594.2858
$ws.Range("L112").Value = 4782.857400000001
$ws.Range("N112").Value = -6998.857400000001
$ws.Range("H137").Value = 4955.5386
$ws.Range("I137").Value = 6113.3213
$ws.Range("K137").Value = 18339.9639
$ws.Range("M137").Value = -15789.9639
$ws.Range("H141").Value = 2871.375
$ws.Range("I141").Value = 1872.5
$ws.Range("J141").Value = 3870.25
$ws.Range("K141").Value = 5617.5
$ws.Range("L141").Value = 11610.75
$ws.Range("M141").Value = -437.5
$ws.Range("N141").Value = -21970.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2923.1277
$ws.Range("I61").Value = 2149.9666
$ws.Range("J61").Value = 4287.5293
$ws.Range("K61").Value = 2149.9666
$ws.Range("L61").Value = 4287.5293
$ws.Range("M61").Value = -1937.9666
$ws.Range("N61").Value = -4711.5293
$ws.Range("H74").Value = 1613.8474
$ws.Range("I74").Value = 959
$ws.Range("K74").Value = 959
$ws.Range("M74").Value = -85
$ws.Range("H77").Value = 1613.8474
$ws.Range("I77").Value = 959
$ws.Range("K77").Value = 4795
$ws.Range("M77").Value = -427
$ws.Range("H82").Value = 38181
$ws.Range("J82").Value = 38181
$ws.Range("L82").Value = 38181
$ws.Range("N82").Value = -38903
$ws.Range("H85").Value = 38181
$ws.Range("J85").Value = 38181
$ws.Range("L85").Value = 38181
$ws.Range("N85").Value = -40677
$ws.Range("H86").Value = 39000
$ws.Range("J86").Value = 39000
$ws.Range("L86").Value = 39000
$ws.Range("N86").Value = -41372
$ws.Range("H89").Value = 39000
$ws.Range("J89").Value = 39000
$ws.Range("L89").Value = 117000
$ws.Range("N89").Value = -128856
$ws.Range("H132").Value = 20316.232
$ws.Range("I132").Value = 24898.205
$ws.Range("K132").Value = 74694.61500000001
$ws.Range("M132").Value = -72164.61500000001
$ws.Range("H136").Value = 2923.1277
$ws.Range("I136").Value = 2149.9666
$ws.Range("J136").Value = 4287.5293
$ws.Range("K136").Value = 6449.899800000001
$ws.Range("L136").Value = 12862.5879
$ws.Range("M136").Value = -3899.899800000001
$ws.Range("N136").Value = -17962.5879
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4682.778
$ws.Range("I86").Value = 7675
$ws.Range("J86").Value = 2289
$ws.Range("K86").Value = 7675
$ws.Range("L86").Value = 2289
$ws.Range("M86").Value = -6552
$ws.Range("N86").Value = -4535
$ws.Range("H89").Value = 4682.778
$ws.Range("I89").Value = 7675
$ws.Range("J89").Value = 2289
$ws.Range("K89").Value = 38375
$ws.Range("L89").Value = 11445
$ws.Range("M89").Value = -32759
$ws.Range("N89").Value = -22677
$ws.Range("H134").Value = 6258.343
$ws.Range("I134").Value = 7349.68
$ws.Range("J134").Value = 3530
$ws.Range("K134").Value = 22049.04
$ws.Range("L134").Value = 10590
$ws.Range("M134").Value = -19514.04
$ws.Range("N134").Value = -15660
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1885.0878
$ws.Range("I31").Value = 1051.5349
$ws.Range("J31").Value = 4445.2856
$ws.Range("K31").Value = 1051.5349
$ws.Range("L31").Value = 4445.2856
$ws.Range("M31").Value = -756.5349000000001
$ws.Range("N31").Value = -5035.2856
$ws.Range("H34").Value = 1885.0878
$ws.Range("I34").Value = 1051.5349
$ws.Range("J34").Value = 4445.2856
$ws.Range("K34").Value = 1051.5349
$ws.Range("L34").Value = 4445.2856
$ws.Range("M34").Value = -849.5349000000001
$ws.Range("N34").Value = -4849.2856
$ws.Range("H58").Value = 2122.7307
$ws.Range("I58").Value = 1511.375
$ws.Range("K58").Value = 1511.375
$ws.Range("M58").Value = -1308.375
$ws.Range("H132").Value = 1983.0435
$ws.Range("I132").Value = 1021.8125
$ws.Range("J132").Value = 4180.143
$ws.Range("K132").Value = 3065.4375
$ws.Range("L132").Value = 12540.429
$ws.Range("M132").Value = -535.4375
$ws.Range("N132").Value = -17600.429
$ws.Range("H134").Value = 1548.9512
$ws.Range("I134").Value = 1039.5454
$ws.Range("K134").Value = 3118.6362
$ws.Range("M134").Value = -583.6361999999999
$ws.Range("H136").Value = 2122.7307
$ws.Range("I136").Value = 1511.375
$ws.Range("K136").Value = 4534.125
$ws.Range("M136").Value = -1984.125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9925
$ws.Range("I56").Value = 9925
$ws.Range("K56").Value = 9925
$ws.Range("M56").Value = -9395
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3471
$ws.Range("I132").Value = 2968.5625
$ws.Range("K132").Value = 8905.6875
$ws.Range("M132").Value = -6375.6875
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1177.8572
$ws.Range("I61").Value = 1191.5385
$ws.Range("K61").Value = 1191.5385
$ws.Range("M61").Value = -989.5385000000001
$ws.Range("H106").Value = 28765.25
$ws.Range("J106").Value = 28765.25
$ws.Range("L106").Value = 28765.25
$ws.Range("N106").Value = -31289.25
$ws.Range("H113").Value = 1177.8572
$ws.Range("I113").Value = 1191.5385
$ws.Range("K113").Value = 1191.5385
$ws.Range("M113").Value = 978.4614999999999
$ws.Range("H122").Value = 13370.889
$ws.Range("I122").Value = 26234.5
$ws.Range("J122").Value = 3080
$ws.Range("K122").Value = 78703.5
$ws.Range("L122").Value = 9240
$ws.Range("M122").Value = -76253.5
$ws.Range("N122").Value = -14140
$ws.Range("H132").Value = 3928.1724
$ws.Range("I132").Value = 3012.8823
$ws.Range("J132").Value = 5224.8335
$ws.Range("K132").Value = 9038.6469
$ws.Range("L132").Value = 15674.5005
$ws.Range("M132").Value = -6508.6469
$ws.Range("N132").Value = -20734.5005
$ws.Range("H136").Value = 3842.2173
$ws.Range("I136").Value = 2491.6428
$ws.Range("J136").Value = 5943.1113
$ws.Range("K136").Value = 7474.928400000001
$ws.Range("L136").Value = 17829.3339
$ws.Range("M136").Value = -4924.928400000001
$ws.Range("N136").Value = -22929.3339
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1694.415
$ws.Range("I132").Value = 917.17145
$ws.Range("J132").Value = 3205.7222
$ws.Range("K132").Value = 2751.51435
$ws.Range("L132").Value = 9617.1666
$ws.Range("M132").Value = -221.5143500000004
$ws.Range("N132").Value = -14677.1666
$ws.Range("H136").Value = 16687255
$ws.Range("I136").Value = 28604628
$ws.Range("J136").Value = 2933.8
$ws.Range("K136").Value = 85813884
$ws.Range("L136").Value = 8801.400000000001
$ws.Range("M136").Value = -85811334
$ws.Range("N136").Value = -13901.4
